# "patentes inseridas no modelo no R"
# Adds the new "Sterman" (I column) values for the existing PeD/patent rows
# (38-45) and appends eleven new parameter rows (46-56) covering the new
# performance / patent-share variables, on the "params" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# --- fill in the "Sterman" (column I) values for the pre-existing PeD rows ---
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(39, 9).Value = 0.001
$ws.Cells.Item(40, 9).Value = 4
$ws.Cells.Item(41, 9).Value = 100000
$ws.Cells.Item(42, 9).Value = 2
$ws.Cells.Item(43, 9).Value = 0.4
$ws.Cells.Item(44, 9).Value = 18
$ws.Cells.Item(45, 9).Value = 10

# --- row 46: aPerfSlope (entered fully: name, friendly name, min/max formulas, unit) ---
$ws.Cells.Item(46, 1).Value = "aPerfSlope"
$ws.Cells.Item(46, 2).Value = "Melhoria em performance por patente que a empresa tem acesso."
$ws.Cells.Item(46, 3).Formula = "=1/30"
$ws.Cells.Item(46, 4).Formula = "=1/30"
$ws.Cells.Item(46, 5).Value = "Unidades de Performance / Patentes"
$ws.Cells.Item(46, 9).Formula = "=1/30"

# --- variable (column A) names typed down for the next few rows first ---
$ws.Cells.Item(47, 1).Value = "aPerfMin"
$ws.Cells.Item(48, 1).Value = "aPerfMax"
$ws.Cells.Item(49, 1).Value = "aSensOfAttractToPerformance"
$ws.Cells.Item(50, 1).Value = "aReferencePerformance"

# --- then the friendly names (column B) were filled back in for Min/Max only ---
$ws.Cells.Item(47, 2).Value = "Índice de Performance Mínimo"
$ws.Cells.Item(48, 2).Value = "Índice de Performance Máximo"

# --- numeric Min/Max/Sens/Reference values ---
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 9).Value = 0

$ws.Cells.Item(48, 3).Value = 10
$ws.Cells.Item(48, 4).Value = 10
$ws.Cells.Item(48, 9).Value = 10

$ws.Cells.Item(49, 3).Value = -4
$ws.Cells.Item(49, 4).Value = -4
$ws.Cells.Item(49, 9).Value = -4

$ws.Cells.Item(50, 3).Value = 10
$ws.Cells.Item(50, 4).Value = 10
$ws.Cells.Item(50, 9).Value = 10

# --- remaining variable (column A) names, rows 51-56 ---
$ws.Cells.Item(51, 1).Value = "aInitialInvestimentoNaoRealizadoPeD"
$ws.Cells.Item(52, 1).Value = "aInitialPatentesRequisitadas"
$ws.Cells.Item(53, 1).Value = "aInitialPatentesEmpresa"
$ws.Cells.Item(54, 1).Value = "aInitialsPatentesEmDominioPublicoUteis"
$ws.Cells.Item(55, 1).Value = "aInitialsInvestimentoPeDDepreciar"
$ws.Cells.Item(56, 1).Value = "aPatentShare"

# --- last friendly name, row 56 ---
$ws.Cells.Item(56, 2).Value = "Share de Patentes"

# --- numeric values for rows 51-56 ---
$ws.Cells.Item(51, 3).Value = 1000
$ws.Cells.Item(51, 4).Value = 1000
$ws.Cells.Item(51, 9).Value = 1000

$ws.Cells.Item(52, 3).Value = 100
$ws.Cells.Item(52, 4).Value = 100
$ws.Cells.Item(52, 9).Value = 100

$ws.Cells.Item(53, 3).Value = 100
$ws.Cells.Item(53, 4).Value = 100
$ws.Cells.Item(53, 9).Value = 100

$ws.Cells.Item(54, 3).Value = 20
$ws.Cells.Item(54, 4).Value = 20
$ws.Cells.Item(54, 9).Value = 20

$ws.Cells.Item(55, 3).Value = 1000000
$ws.Cells.Item(55, 4).Value = 1000000
$ws.Cells.Item(55, 9).Value = 1000000

$ws.Cells.Item(56, 3).Value = 0.5
$ws.Cells.Item(56, 4).Value = 0.5
$ws.Cells.Item(56, 9).Value = 0.5

# --- column width tweaks (approximate to the nearest width this engine can store) ---
$ws.Columns.Item(1).ColumnWidth = 37.8333333333333
$ws.Columns.Item(3).ColumnWidth = 12.1666666666667
$ws.Columns.Item(4).ColumnWidth = 9.66666666666667
$ws.Columns.Item(9).ColumnWidth = 7.83333333333333

# --- restore the view/selection close to what the authored workbook shows ---
$ws.Range("E60").Select()

Write-Output "patentes inseridas no modelo no R - edit applied"
